# Generate Report for Handback
# Adds a second handback entry (ee2b127f-e981-41b1-82e4-fb07c0b804a4.md) next to
# the existing one (renamed from 1f029654-... to b84f5bbe-...) across all three
# sheets (Overview, zh-cn, de-de), and extends each table by one row.

$wb = $excel.ActiveWorkbook

$oldGuid = "1f029654-b100-467c-9a79-28e6546fe1bc"
$newGuid1 = "b84f5bbe-6c67-4ae1-b690-2f3d283370b5"
$newGuid2 = "ee2b127f-e981-41b1-82e4-fb07c0b804a4"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsO = $wb.Worksheets.Item("Overview")

# Rename the existing handback file (row 2) and refresh its timestamp.
$wsO.Range("A2").Value = "$newGuid1.md"
$wsO.Range("B2").Value = "e2e\$newGuid1.md"
$wsO.Range("G2").Value = "2016-08-15 09:15:12"

# Update the existing hyperlink's display text to match (target URL unchanged).
$wsO.Range("B2").Hyperlinks.Delete()
$wsO.Hyperlinks.Add($wsO.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bce4fc9b9dab65c68a62f0fd68492749f6ae58d4/e2e/$newGuid1.md", [type]::Missing, [type]::Missing, "e2e\$newGuid1.md") | Out-Null
$wsO.Range("B2").Style = "HyperLink"

# Append a new row to the Overview table for the second handback file.
$loO = $wsO.ListObjects.Item(1)
$loO.ListRows.Add() | Out-Null

$wsO.Range("A3").Value = "$newGuid2.md"
$wsO.Range("B3").Value = "e2e\$newGuid2.md"
$wsO.Range("C3").Value = ".md"
$wsO.Range("E3").Value = "Handed back: in sync with en-US"
$wsO.Range("F3").Value = "Handed back: in sync with en-US"
$wsO.Range("G3").Value = "2016-08-15 09:15:12"
$wsO.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsO.Hyperlinks.Add($wsO.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bce4fc9b9dab65c68a62f0fd68492749f6ae58d4/e2e/$newGuid2.md", [type]::Missing, [type]::Missing, "e2e\$newGuid2.md") | Out-Null
$wsO.Range("B3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZ = $wb.Worksheets.Item("zh-cn")

# Update row 2 (existing file renamed + new xliff + new timestamps).
$wsZ.Range("A2").Value = "$newGuid1.md"
$wsZ.Range("G2").Value = "$newGuid1.d02e7ee50af643d5230de531e41d00ef1bca9c60.zh-cn.xlf"
$wsZ.Range("H2").Value = "2016-08-15 09:14:59"
$wsZ.Range("I2").Value = "$newGuid1.md"
$wsZ.Range("J2").Value = "$newGuid1.d02e7ee50af643d5230de531e41d00ef1bca9c60.zh-cn.xlf"
$wsZ.Range("K2").Value = "2016-08-15 09:15:28"

$wsZ.Range("A2").Hyperlinks.Delete()
$wsZ.Hyperlinks.Add($wsZ.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bce4fc9b9dab65c68a62f0fd68492749f6ae58d4/e2e/$newGuid1.md", [type]::Missing, [type]::Missing, "$newGuid1.md") | Out-Null
$wsZ.Range("A2").Style = "HyperLink"

$wsZ.Range("I2").Hyperlinks.Delete()
$wsZ.Hyperlinks.Add($wsZ.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b1a02ff2ecd3a8ac1d829ea9dbec1bb7258b114e/e2e/$newGuid1.md", [type]::Missing, [type]::Missing, "$newGuid1.md") | Out-Null
$wsZ.Range("I2").Style = "HyperLink"

# Append row 3 for the second handback file.
$loZ = $wsZ.ListObjects.Item(1)
$loZ.ListRows.Add() | Out-Null

$wsZ.Range("A3").Value = "$newGuid2.md"
$wsZ.Range("B3").Value = ".md"
$wsZ.Range("C3").Value = "Handed back: in sync with en-US"
$wsZ.Range("D3").Value = "e2e"
$wsZ.Range("E3").Value = "ht"
$wsZ.Range("F3").Value = "True"
$wsZ.Range("G3").Value = "$newGuid2.3f8b2a439c624253ccbc80bf855779a74d8b9e17.zh-cn.xlf"
$wsZ.Range("H3").Value = "2016-08-15 09:14:59"
$wsZ.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZ.Range("I3").Value = "$newGuid2.md"
$wsZ.Range("J3").Value = "$newGuid2.3f8b2a439c624253ccbc80bf855779a74d8b9e17.zh-cn.xlf"
$wsZ.Range("K3").Value = "2016-08-15 09:15:28"
$wsZ.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZ.Range("L3").Value = ""
$wsZ.Range("M3").Value = "True"
$wsZ.Range("N3").Value = ""
$wsZ.Range("O3").Value = "False"
$wsZ.Range("P3").Value = ""

$wsZ.Hyperlinks.Add($wsZ.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bce4fc9b9dab65c68a62f0fd68492749f6ae58d4/e2e/$newGuid2.md", [type]::Missing, [type]::Missing, "$newGuid2.md") | Out-Null
$wsZ.Range("A3").Style = "HyperLink"
$wsZ.Hyperlinks.Add($wsZ.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b1a02ff2ecd3a8ac1d829ea9dbec1bb7258b114e/e2e/$newGuid2.md", [type]::Missing, [type]::Missing, "$newGuid2.md") | Out-Null
$wsZ.Range("I3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsD = $wb.Worksheets.Item("de-de")

# Update row 2 (existing file renamed + new xliff + new timestamps).
$wsD.Range("A2").Value = "$newGuid1.md"
$wsD.Range("G2").Value = "$newGuid1.d02e7ee50af643d5230de531e41d00ef1bca9c60.de-de.xlf"
$wsD.Range("H2").Value = "2016-08-15 09:15:12"
$wsD.Range("I2").Value = "$newGuid1.md"
$wsD.Range("J2").Value = "$newGuid1.d02e7ee50af643d5230de531e41d00ef1bca9c60.de-de.xlf"
$wsD.Range("K2").Value = "2016-08-15 09:15:35"

$wsD.Range("A2").Hyperlinks.Delete()
$wsD.Hyperlinks.Add($wsD.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bce4fc9b9dab65c68a62f0fd68492749f6ae58d4/e2e/$newGuid1.md", [type]::Missing, [type]::Missing, "$newGuid1.md") | Out-Null
$wsD.Range("A2").Style = "HyperLink"

$wsD.Range("I2").Hyperlinks.Delete()
$wsD.Hyperlinks.Add($wsD.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/ae90bf7a8bdefd76620b17ad4dd6726696a323e7/e2e/$newGuid1.md", [type]::Missing, [type]::Missing, "$newGuid1.md") | Out-Null
$wsD.Range("I2").Style = "HyperLink"

# Append row 3 for the second handback file.
$loD = $wsD.ListObjects.Item(1)
$loD.ListRows.Add() | Out-Null

$wsD.Range("A3").Value = "$newGuid2.md"
$wsD.Range("B3").Value = ".md"
$wsD.Range("C3").Value = "Handed back: in sync with en-US"
$wsD.Range("D3").Value = "e2e"
$wsD.Range("E3").Value = "ht"
$wsD.Range("F3").Value = "True"
$wsD.Range("G3").Value = "$newGuid2.3f8b2a439c624253ccbc80bf855779a74d8b9e17.de-de.xlf"
$wsD.Range("H3").Value = "2016-08-15 09:15:12"
$wsD.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsD.Range("I3").Value = "$newGuid2.md"
$wsD.Range("J3").Value = "$newGuid2.3f8b2a439c624253ccbc80bf855779a74d8b9e17.de-de.xlf"
$wsD.Range("K3").Value = "2016-08-15 09:15:35"
$wsD.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsD.Range("L3").Value = ""
$wsD.Range("M3").Value = "True"
$wsD.Range("N3").Value = ""
$wsD.Range("O3").Value = "False"
$wsD.Range("P3").Value = ""

$wsD.Hyperlinks.Add($wsD.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bce4fc9b9dab65c68a62f0fd68492749f6ae58d4/e2e/$newGuid2.md", [type]::Missing, [type]::Missing, "$newGuid2.md") | Out-Null
$wsD.Range("A3").Style = "HyperLink"
$wsD.Hyperlinks.Add($wsD.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/ae90bf7a8bdefd76620b17ad4dd6726696a323e7/e2e/$newGuid2.md", [type]::Missing, [type]::Missing, "$newGuid2.md") | Out-Null
$wsD.Range("I3").Style = "HyperLink"

$wb.Save()
